# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Apio" (Terminal Hortofrutícola Agro
# Chillán) as row 116, shifting the existing rows 116-142 down to 117-143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 116:142 down by inserting a new blank row at 116.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range('A116').Value = 7
$ws.Range('B116').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C116').Value = 'Ñuble'
$ws.Range('D116').Value = 44508
$ws.Range('E116').Value = 16
$ws.Range('F116').Value = 100112017
$ws.Range('G116').Value = 'Apio'
$ws.Range('H116').Value = 'Americana (o)'
$ws.Range('I116').Value = 'Primera'
$ws.Range('J116').Value = 120
$ws.Range('K116').Value = 8000
$ws.Range('L116').Value = 9000
$ws.Range('M116').Value = 8500
$ws.Range('N116').Value = '$/docena de matas'
$ws.Range('O116').Value = 'Provincia del Elquí'
$ws.Range('P116').Value = 1417
$ws.Range('Q116').Value = 6
$ws.Range('R116').Value = 'Hortaliza'
